$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the data range as text (string) cells, matching the original
# inline-string typing used throughout this table.
$ws.Range("B2:K21").NumberFormat = "@"

# Row 2: Real Madrid
$ws.Range("B2").Value = "1"
$ws.Range("C2").Value = " Real Madrid"
$ws.Range("D2").Value = "20"
$ws.Range("E2").Value = "14"
$ws.Range("F2").Value = "4"
$ws.Range("G2").Value = "2"
$ws.Range("H2").Value = "41"
$ws.Range("I2").Value = "17"
$ws.Range("J2").Value = "+24"
$ws.Range("K2").Value = "46"

# Row 3: Sevilla
$ws.Range("B3").Value = "2"
$ws.Range("C3").Value = " Sevilla"
$ws.Range("D3").Value = "19"
$ws.Range("E3").Value = "12"
$ws.Range("F3").Value = "5"
$ws.Range("G3").Value = "2"
$ws.Range("H3").Value = "30"
$ws.Range("I3").Value = "13"
$ws.Range("J3").Value = "+17"
$ws.Range("K3").Value = "41"

# Row 4: Betis
$ws.Range("B4").Value = "3"
$ws.Range("C4").Value = " Betis"
$ws.Range("D4").Value = "19"
$ws.Range("E4").Value = "10"
$ws.Range("F4").Value = "3"
$ws.Range("G4").Value = "6"
$ws.Range("H4").Value = "32"
$ws.Range("I4").Value = "23"
$ws.Range("J4").Value = "+9"
$ws.Range("K4").Value = "33"

# Row 5: Atlético Madrid
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = " Atlético Madrid"
$ws.Range("D5").Value = "19"
$ws.Range("E5").Value = "9"
$ws.Range("F5").Value = "5"
$ws.Range("G5").Value = "5"
$ws.Range("H5").Value = "31"
$ws.Range("I5").Value = "22"
$ws.Range("J5").Value = "+9"
$ws.Range("K5").Value = "32"

# Row 6: Barcelona
$ws.Range("B6").Value = "5"
$ws.Range("C6").Value = " Barcelona"
$ws.Range("D6").Value = "19"
$ws.Range("E6").Value = "8"
$ws.Range("F6").Value = "7"
$ws.Range("G6").Value = "4"
$ws.Range("H6").Value = "30"
$ws.Range("I6").Value = "22"
$ws.Range("J6").Value = "+8"
$ws.Range("K6").Value = "31"

# Row 7: Rayo Vallecano
$ws.Range("B7").Value = "6"
$ws.Range("C7").Value = " Rayo Vallecano"
$ws.Range("D7").Value = "19"
$ws.Range("E7").Value = "9"
$ws.Range("F7").Value = "3"
$ws.Range("G7").Value = "7"
$ws.Range("H7").Value = "26"
$ws.Range("I7").Value = "20"
$ws.Range("J7").Value = "+6"
$ws.Range("K7").Value = "30"

# Row 8: Real Sociedad
$ws.Range("B8").Value = "7"
$ws.Range("C8").Value = " Real Sociedad"
$ws.Range("D8").Value = "19"
$ws.Range("E8").Value = "8"
$ws.Range("F8").Value = "6"
$ws.Range("G8").Value = "5"
$ws.Range("H8").Value = "21"
$ws.Range("I8").Value = "21"
$ws.Range("J8").Value = "0"
$ws.Range("K8").Value = "30"

# Row 9: Villarreal
$ws.Range("B9").Value = "8"
$ws.Range("C9").Value = " Villarreal"
$ws.Range("D9").Value = "19"
$ws.Range("E9").Value = "7"
$ws.Range("F9").Value = "7"
$ws.Range("G9").Value = "5"
$ws.Range("H9").Value = "31"
$ws.Range("I9").Value = "20"
$ws.Range("J9").Value = "+11"
$ws.Range("K9").Value = "28"

# Row 10: Valencia
$ws.Range("B10").Value = "9"
$ws.Range("C10").Value = " Valencia"
$ws.Range("D10").Value = "19"
$ws.Range("E10").Value = "7"
$ws.Range("F10").Value = "7"
$ws.Range("G10").Value = "5"
$ws.Range("H10").Value = "31"
$ws.Range("I10").Value = "28"
$ws.Range("J10").Value = "+3"
$ws.Range("K10").Value = "28"

# Row 11: Athletic Club
$ws.Range("B11").Value = "10"
$ws.Range("C11").Value = " Athletic Club"
$ws.Range("D11").Value = "20"
$ws.Range("E11").Value = "6"
$ws.Range("F11").Value = "9"
$ws.Range("G11").Value = "5"
$ws.Range("H11").Value = "20"
$ws.Range("I11").Value = "17"
$ws.Range("J11").Value = "+3"
$ws.Range("K11").Value = "27"

# Row 12: Espanyol
$ws.Range("B12").Value = "11"
$ws.Range("C12").Value = " Espanyol"
$ws.Range("D12").Value = "19"
$ws.Range("E12").Value = "7"
$ws.Range("F12").Value = "5"
$ws.Range("G12").Value = "7"
$ws.Range("H12").Value = "22"
$ws.Range("I12").Value = "22"
$ws.Range("J12").Value = "0"
$ws.Range("K12").Value = "26"

# Row 13: Celta Vigo
$ws.Range("B13").Value = "12"
$ws.Range("C13").Value = " Celta Vigo"
$ws.Range("D13").Value = "19"
$ws.Range("E13").Value = "6"
$ws.Range("F13").Value = "5"
$ws.Range("G13").Value = "8"
$ws.Range("H13").Value = "22"
$ws.Range("I13").Value = "22"
$ws.Range("J13").Value = "0"
$ws.Range("K13").Value = "23"

# Row 14: Granada
$ws.Range("B14").Value = "13"
$ws.Range("C14").Value = " Granada"
$ws.Range("D14").Value = "19"
$ws.Range("E14").Value = "5"
$ws.Range("F14").Value = "8"
$ws.Range("G14").Value = "6"
$ws.Range("H14").Value = "23"
$ws.Range("I14").Value = "26"
$ws.Range("J14").Value = "-3"
$ws.Range("K14").Value = "23"

# Row 15: Osasuna
$ws.Range("B15").Value = "14"
$ws.Range("C15").Value = " Osasuna"
$ws.Range("D15").Value = "19"
$ws.Range("E15").Value = "5"
$ws.Range("F15").Value = "7"
$ws.Range("G15").Value = "7"
$ws.Range("H15").Value = "18"
$ws.Range("I15").Value = "25"
$ws.Range("J15").Value = "-7"
$ws.Range("K15").Value = "22"

# Row 16: Mallorca
$ws.Range("B16").Value = "15"
$ws.Range("C16").Value = " Mallorca"
$ws.Range("D16").Value = "19"
$ws.Range("E16").Value = "4"
$ws.Range("F16").Value = "8"
$ws.Range("G16").Value = "7"
$ws.Range("H16").Value = "17"
$ws.Range("I16").Value = "28"
$ws.Range("J16").Value = "-11"
$ws.Range("K16").Value = "20"

# Row 17: Getafe
$ws.Range("B17").Value = "16"
$ws.Range("C17").Value = " Getafe"
$ws.Range("D17").Value = "19"
$ws.Range("E17").Value = "4"
$ws.Range("F17").Value = "6"
$ws.Range("G17").Value = "9"
$ws.Range("H17").Value = "13"
$ws.Range("I17").Value = "20"
$ws.Range("J17").Value = "-7"
$ws.Range("K17").Value = "18"

# Row 18: Elche
$ws.Range("B18").Value = "17"
$ws.Range("C18").Value = " Elche"
$ws.Range("D18").Value = "19"
$ws.Range("E18").Value = "3"
$ws.Range("F18").Value = "7"
$ws.Range("G18").Value = "9"
$ws.Range("H18").Value = "18"
$ws.Range("I18").Value = "27"
$ws.Range("J18").Value = "-9"
$ws.Range("K18").Value = "16"

# Row 19: Alavés
$ws.Range("B19").Value = "18"
$ws.Range("C19").Value = " Alavés"
$ws.Range("D19").Value = "19"
$ws.Range("E19").Value = "4"
$ws.Range("F19").Value = "4"
$ws.Range("G19").Value = "11"
$ws.Range("H19").Value = "16"
$ws.Range("I19").Value = "30"
$ws.Range("J19").Value = "-14"
$ws.Range("K19").Value = "16"

# Row 20: Cádiz
$ws.Range("B20").Value = "19"
$ws.Range("C20").Value = " Cádiz"
$ws.Range("D20").Value = "19"
$ws.Range("E20").Value = "2"
$ws.Range("F20").Value = "8"
$ws.Range("G20").Value = "9"
$ws.Range("H20").Value = "15"
$ws.Range("I20").Value = "32"
$ws.Range("J20").Value = "-17"
$ws.Range("K20").Value = "14"

# Row 21: Levante
$ws.Range("B21").Value = "20"
$ws.Range("C21").Value = " Levante"
$ws.Range("D21").Value = "19"
$ws.Range("E21").Value = "0"
$ws.Range("F21").Value = "8"
$ws.Range("G21").Value = "11"
$ws.Range("H21").Value = "19"
$ws.Range("I21").Value = "41"
$ws.Range("J21").Value = "-22"
$ws.Range("K21").Value = "8"

